# Form the consolidated report: fill in the "Absent" column (H) so that
# every row has a value (Absent = 1 when Real attendance is 0, else 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Absent value needs to flip from 0 to 1 (student was absent).
$absentRows = @(10, 13, 19)
foreach ($r in $absentRows) {
    $ws.Cells.Item($r, 8).Value = 1
}

# Rows whose Absent cell was blank and needs to be explicitly 0 (student was present).
$presentRows = @(11, 14, 20)
foreach ($r in $presentRows) {
    $ws.Cells.Item($r, 8).Value = 0
}
